$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "OBSERVATIONS" column headers to more specific labels
$ws.Range("A1").Value = "OBSERVATIONS N UMAP FIT"
$ws.Range("D1").Value = "OBSERVATIONS N UMAP TRANSFORM"
$ws.Range("G1").Value = "OBSERVATIONS N HDBSCAN FIT"
$ws.Range("J1").Value = "OBSERVATIONS N HDBSCAN TRANSFORM"
$ws.Range("M1").Value = "OBSERVATIONS N XGBOOST FIT"
$ws.Range("P1").Value = "OBSERVATIONS N XGBOOST TRANSFORM"
$ws.Range("Q1").Value = "XGBOOST TRANSFORM GPU"
$ws.Range("R1").Value = "XGBOOST TRANSFORM CPU"

# Resize columns to fit the new (longer) header text (mirrors Excel's own
# "best fit" recalculation of column width after the header text changed)
$ws.Columns.Item(1).ColumnWidth = 24.330729166666668
$ws.Columns.Item(4).ColumnWidth = 32.998697916666664
$ws.Columns.Item(7).ColumnWidth = 27.666666666666668
$ws.Columns.Item(10).ColumnWidth = 36.330729166666664
$ws.Columns.Item(13).ColumnWidth = 27.666666666666668
$ws.Columns.Item(16).ColumnWidth = 36.330729166666664
$ws.Columns.Item(17).ColumnWidth = 24.498697916666668
$ws.Columns.Item(18).ColumnWidth = 24.166666666666668

# Update view: zoom and selection
$ws.Range("B7").Select()
$excel.ActiveWindow.Zoom = 125
